$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Il12b"
$ws.Range("C2").Value = "Il23r"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.091959
$ws.Range("H2").Value = 0.275877
$ws.Range("I2").Value = 0.04433703455491324
$ws.Range("J2").Value = 0.04433703455491323
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.014443
$ws.Range("N2").Value = 0.043329
$ws.Range("O2").Value = 0.7689945869198687
$ws.Range("P2").Value = 0.7689945869198687
$ws.Range("Q2").Value = 0.001328163837
$ws.Range("R2").Value = 0.011953474533
$ws.Range("S2").Value = 0.03409493957280745
$ws.Range("T2").Value = 0.03409493957280745

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il12b"
$ws.Range("C3").Value = "Il23r"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.091959
$ws.Range("H3").Value = 0.275877
$ws.Range("I3").Value = 0.04433703455491324
$ws.Range("J3").Value = 0.04433703455491323
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.004338666666666667
$ws.Range("N3").Value = 0.013016
$ws.Range("O3").Value = 0.2310054130801314
$ws.Range("P3").Value = 0.2310054130801313
$ws.Range("Q3").Value = 0.000398979448
$ws.Range("R3").Value = 0.003590815032
$ws.Range("S3").Value = 0.01024209498210579
$ws.Range("T3").Value = 0.01024209498210579

# Row 4
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Il12b"
$ws.Range("C4").Value = "Il23r"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.982131
$ws.Range("H4").Value = 5.946393
$ws.Range("I4").Value = 0.9556629654450868
$ws.Range("J4").Value = 0.9556629654450867
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.014443
$ws.Range("N4").Value = 0.043329
$ws.Range("O4").Value = 0.7689945869198687
$ws.Range("P4").Value = 0.7689945869198687
$ws.Range("Q4").Value = 0.028627918033
$ws.Range("R4").Value = 0.257651262297
$ws.Range("S4").Value = 0.7348996473470613
$ws.Range("T4").Value = 0.7348996473470611

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Il12b"
$ws.Range("C5").Value = "Il23r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.982131
$ws.Range("H5").Value = 5.946393
$ws.Range("I5").Value = 0.9556629654450868
$ws.Range("J5").Value = 0.9556629654450867
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.004338666666666667
$ws.Range("N5").Value = 0.013016
$ws.Range("O5").Value = 0.2310054130801314
$ws.Range("P5").Value = 0.2310054130801313
$ws.Range("Q5").Value = 0.008599805698666668
$ws.Range("R5").Value = 0.077398251288
$ws.Range("S5").Value = 0.2207633180980256
$ws.Range("T5").Value = 0.2207633180980255
